$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A3").Value = 112492987
$ws.Range("B3").Value = 90818
$ws.Range("E3").Value = 4368
$ws.Range("Q3").Value = 524928
$ws.Range("R3").Value = 6540273
$ws.Range("S3").Value = 25

# Text cells
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Dofttaggsvamp"
$ws.Range("G3").Value = "Hydnellum suaveolens"
$ws.Range("H3").Value = "(Scop.:Fr.) P. Karst."
$ws.Range("P3").Value = "Kattfall NO-ut, Nrk"
$ws.Range("T3").Value = "Örebro"
$ws.Range("U3").Value = "Hallsberg"
$ws.Range("V3").Value = "Närke"
$ws.Range("W3").Value = "Svennevad"

# Date-looking text cells - force text format so they are stored as plain
# strings (matching source data) instead of being coerced to date serials.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-10-01"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-10-01"

# Boolean cells
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# Trailing text cells
$ws.Range("AW3").Value = "Arne Holmer"
$ws.Range("AX3").Value = "Arne Holmer"

# Empty placeholder cells present in the row (no value, default style) -
# touching a formatting no-op materializes the cell without giving it
# content, matching the source row's empty inline-string cells.
foreach ($addr in @("I3", "J3", "K3", "N3", "AF3", "AT3", "AY3")) {
    $ws.Range($addr).Font.Bold = $false
}
